$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header changes ---
# Column W ("additional") is renamed to "Non DPR"
$ws.Range("W1").Value = "Non DPR"

# New column X ("Non Surveyed") is appended, reusing W1 header formatting
$ws.Range("W1").Copy() | Out-Null
$ws.Range("X1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("X1").Value = "Non Surveyed"

# --- Data changes (columns B:X, rows 2:11) ---
# Each inner array holds the values for columns B..X (acsr .. Non Surveyed)
# for one row (BISHNUPUR, CHANDEL, ..., UKHRUL, Sum).
$data = @(
    @(14.46, 23.47, 19.68, 1331, 20, 3, 3, 4, 6.35, 11.8125, 84, 22.11, 18.405, 1171, 2, 1, 3, 49, 0, 50, 54, 0, 4),  # BISHNUPUR
    @(26.6, 639.3700000000002, 89.37999999999992, 3373, 96, 0, 3, 2, 4.720000000000001, 11.025, 82, 31.97000000000001, 71.73000000000002, 2601, 0, 4, 7, 179, 24, 406, 412, 2, 4),  # CHANDEL
    @(89.425, 48.14599999999998, 66.38399999999999, 4302, 206, 0, 4, 99, 43.474, 275.5505000000001, 1753, 48.00899999999999, 66.14799999999998, 2751, 0, 4, 99, 127, 69, 173, 164, 41, 8),  # CHURACHANDPUR
    @(8.07, 6.420000000000001, 8.37, 710, 44, 0, 11, 6, 10.54, 23.9085, 238, 9.16, 9.69, 435, 0, 11, 8, 30, 3, 55, 57, 3, 5),  # IMPHAL EAST
    @(10.7, 12.51, 15.655, 1004, 48, 0, 18, 6, 8.33, 14.829, 214, 13.51, 15.875, 916, 1, 14, 4, 31, 0, 56, 63, 13, 7),  # IMPHAL WEST
    @(11.324, 43.62299999999998, 82.642, 3320, 46, 0, 8, 15, 4.824999999999999, 3.16575, 72, 12.761, 20.177, 1933, 0, 3, 8, 50, 0, 406, 413, 111, 0),  # SENAPATI
    @(64.1, 71.92999999999996, 51.19999999999999, 3139, 53, 0, 9, 12, 2.6, 0, 33, 30.55000000000001, 25.7, 1312, 0, 2, 1, 37, 0, 130, 135, 0, 2),  # TAMENGLONG
    @(50.343, 23.8555, 36.25400000000001, 1629, 48, 0, 6, 19, 10.173, 26.9073, 156, 20.731, 25.64899999999999, 1151, 0, 7, 21, 60, 19, 78, 137, 13, 6),  # THOUBAL
    @(178.575, 45.60999999999996, 33.67899999999999, 1778, 66, 0, 13, 21, 30.31, 60.60600000000001, 379, 24.9, 18.959, 1029, 0, 19, 29, 70, 35, 124, 123, 54, 0),  # UKHRUL
    @(453.597, 914.9345000000001, 403.2439999999999, 20586, 627, 3, 75, 184, 121.322, 427.8045500000001, 3011, 213.701, 272.333, 13299, 3, 65, 180, 633, 150, 1478, 1558, 237, 36)  # ∑
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}
